$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller - testing")
$ws.Activate()

# --- Row 7 (patients_controller): index view missing, destroy method present ---
$ws.Cells.Item(7,1).Interior.Color = 255
$ws.Cells.Item(7,2).Value = "yes"
$ws.Cells.Item(7,4).Value = "no"
$ws.Cells.Item(7,5).Value = "you are missing the index view, so at this point in time all test which require an index fail. Good news this controller has a destroy method"
$ws.Rows.Item(7).RowHeight = 105

# --- Row 8 (admins_controller): testing commenced, lots of bugs found ---
$ws.Cells.Item(8,1).Interior.Color = 255
$ws.Cells.Item(8,2).Value = "yes"
$ws.Cells.Item(8,4).Value = ":show"
$ws.Cells.Item(8,5).Value = ":index does not work (returns nil), :new does not work returns nil, :edit doesnt work (does not redirect and returns nil)"
$ws.Cells.Item(8,6).VerticalAlignment = -4108
$ws.Cells.Item(8,6).Value = "Leanne check :create again"
$ws.Rows.Item(8).RowHeight = 90

# --- Row 9: follow-up bug note on :update ---
$ws.Cells.Item(9,5).Value = ":update remoces password, when no changes are suppose to occur"
$ws.Rows.Item(9).RowHeight = 60

[void]$ws.Range("E10").Select()

Write-Host "done"
